# Reecriture de la classe décor
# Insert a new task row ("Remettre le code aux normes" / "Respecter les
# règles de nommage d'après le document") above the existing "Reflexion sur
# la gestion de Map" row, shifting every row from 7 down to 16 one row
# lower (7->8, 8->9, ... 16->17). Also bump two progress percentages and
# refresh the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Shift rows 16..7 down into rows 17..8, copying both formatting
#        and values so style indices / shared strings line up with the
#        pre-existing rows (processed bottom-up so sources aren't
#        clobbered before they are read). ---
for ($r = 16; $r -ge 7; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("A" + $srcRow + ":F" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":F" + $dstRow)
    $src.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
    $dst.Value2 = $src.Value2
}

# --- 2) Populate the now-vacated row 7 with the new task. Pull formats
#        from cells that already carry the desired style so no new
#        cellXfs entries get created. ---
$ws.Range("A2:D2").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E3").Copy()
$ws.Range("E7").PasteSpecial(-4122)      # xlPasteFormats (percentage style)

$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial(-4122)      # xlPasteFormats (orange, bold font)

$ws.Range("A7").Value2 = "Remettre le code aux normes"
$ws.Range("B7").Value2 = "Respecter les règles de nommage d'après le document"
$ws.Range("C7").Value2 = "Jo / Pizzi "
$ws.Range("D7").Value2 = "En attente "
$ws.Range("E7").Value2 = 0
$ws.Range("F7").Value2 = "normal"

# --- 3) Progress percentage tweaks. ---
$ws.Range("E3").Value2 = 0.4
$ws.Range("E5").Value2 = 0.35

# --- 4) Refresh the selection / scroll position. ---
$ws.Range("E5").Select() | Out-Null
